# Applies the "Api" worksheet addition to the workbook, matching the
# commit "changes in pom and apitestscenarios".
#
# Summary of behaviour implemented:
#  - Insert a new worksheet named "Api" right after "Sheet1" (sheetId 2).
#  - Populate it with a small scenario/company id lookup table.
#  - Format the populated cells as Text (so numeric-looking ids like
#    "6044" stay text rather than becoming numbers).
#  - Size the columns to fit their content.
#  - Move the "active"/selected tab from Sheet1 to the new Api sheet, and
#    update the remembered selection on both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Sheet1 keeps its previous selection cell logic but the commit moved the
# remembered selection from D4 to B4.
[void]$ws1.Range("B4").Select()

# Insert the new sheet directly after Sheet1.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Api"

# Pre-format the cells that will hold the id values as Text *before*
# entering the values, so that "6044" / "6045" / "1739" are stored as
# text rather than being auto-converted to numbers.
$ws2.Range("B1:C3").NumberFormat = "@"
$ws2.Range("A2").NumberFormat = "@"

# Fill in the values (order matches how the shared strings table was
# originally built).
$ws2.Range("A2").Value = "API Data01"
$ws2.Range("B1").Value = "ScenarioID"
$ws2.Range("C1").Value = "CompanyID"
$ws2.Range("B2").Value = "6044"
$ws2.Range("B3").Value = "6045"
$ws2.Range("C2").Value = "1739"
$ws2.Range("C3").Value = "1739"

# Size the columns to fit their contents (bestFit-style autosizing).
$ws2.Columns("A").ColumnWidth = 10.61
$ws2.Columns("B").ColumnWidth = 19.17
$ws2.Columns("C").ColumnWidth = 32.83
$ws2.Columns("D").ColumnWidth = 8.83

# The new sheet becomes the active/selected tab, with D1 as the
# remembered selection.
[void]$ws2.Range("D1").Select()

# Match the printed page orientation recorded for the new sheet.
$ws2.PageSetup.Orientation = 1

Write-Output "Added Api worksheet"
